# The "Artfynd" sheet holds one observation record per data row (rows 2-4,
# header in row 1). This edit rotates the three records: the record that
# used to be on row 4 moves to row 2, the record on row 2 moves to row 3,
# and the record on row 3 moves to row 4. Concretely:
#     new row 2  <=  old row 4
#     new row 3  <=  old row 2
#     new row 4  <=  old row 3
#
# Columns differ slightly per record (e.g. only the lichen record that
# used to be on row 4 has an "Auktor"/H value, only the lycopod records
# that used to be on rows 2/3 have values in J/K/L, etc.), so some cells
# gain content, some lose it, and some become an empty (but still
# present) text cell, while plenty of columns (Kommun, Provins, the
# "Ej aterfunnen" flags, ...) happen to hold identical values across all
# three records and must be left completely untouched.
#
# Strategy:
#  1. Snapshot every relevant cell of rows 2-4 first (so the rotation
#     doesn't clobber source data while writing).
#  2. For every destination cell, compare what it should become (the
#     rotated source value) with what it currently holds, and only touch
#     cells whose content actually needs to change -- this avoids
#     incidental format/style churn (e.g. re-assigning an unchanged
#     string would otherwise still rewrite the cell) and keeps untouched
#     columns byte-for-byte identical.
#  3. Numeric columns get a plain numeric .Value. Text columns get their
#     value written with a leading apostrophe so Excel always stores it
#     as literal text -- this prevents Excel from reinterpreting
#     date-like content (e.g. "2016-10-02") as a real date/number, and
#     guarantees an empty string is stored as a real (present) empty
#     text cell instead of silently clearing/removing the cell.
#  4. A destination cell whose source does not exist at all (that column
#     was never populated for the source record) is cleared so it ends
#     up absent too, matching the source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric columns (stored as t="n" in the sheet).
$numericColumns = @("A", "B", "E", "Q", "R", "S")

# Every other populated column across rows 2-4 (text columns; the
# boolean flag columns AD/AE/AG never actually change between records
# here, so they are intentionally excluded and left untouched).
$textColumns = @("C", "D", "F", "G", "H", "I", "J", "K", "L", "P", "T", "U", "V", "W", `
                  "Y", "Z", "AA", "AB", "AC", "AI", "AR", "AT", "AW", "AX", "AY")

$allColumns = $numericColumns + $textColumns

function Get-RowSnapshot([int]$rowIndex) {
    $snap = @{}
    foreach ($col in $allColumns) {
        $snap[$col] = $ws.Range($col + $rowIndex).Value2
    }
    return $snap
}

# Snapshot the current (pre-edit) contents of rows 2, 3 and 4 up front so
# that writing the new row 2 doesn't affect the data we still need to read
# for the new row 3/4 (and vice versa).
$oldRows = @{
    2 = Get-RowSnapshot 2
    3 = Get-RowSnapshot 3
    4 = Get-RowSnapshot 4
}

# new row index -> old row index it should inherit data from.
$mapping = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($newRowIndex in @(2, 3, 4)) {
    $srcRow = $oldRows[$mapping[$newRowIndex]]
    $curRow = $oldRows[$newRowIndex]

    foreach ($col in $allColumns) {
        $newValue = $srcRow[$col]
        $curValue = $curRow[$col]

        if ($newValue -eq $curValue) {
            # Nothing to do: the value at this position is already
            # correct (either both empty/absent, or already identical).
            continue
        }

        $cell = $ws.Range($col + $newRowIndex)

        if ($null -eq $newValue) {
            # Source record never had this column populated -> make sure
            # the destination doesn't have it either.
            $cell.ClearContents() | Out-Null
        }
        elseif ($numericColumns -contains $col) {
            $cell.Value = $newValue
        }
        else {
            # Leading apostrophe forces literal text: prevents Excel from
            # re-parsing dates/numbers, and guarantees an empty string is
            # stored as a real empty text cell rather than clearing it.
            $cell.Value = "'" + $newValue
        }
    }
}
